# Supplementary Table 5: update "Secondary" age-group row (row 10) values.
# Values are stored as text strings (inlineStr), not numbers, so we force
# the Text number format before writing, then reset the style back to
# Normal so no residual style index is left on the cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B10" = "0.21"
    "C10" = "0.42"
    "D10" = "0.3"
    "E10" = "0.45"
    "F10" = "0.39"
    "G10" = "0.53"
    "H10" = "0.47"
    "I10" = "0.53"
    "J10" = "0.56"
    "K10" = "0.62"
    "L10" = "0.59"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
